# Rename three of the "syst*_c" systematic-uncertainty column headers to
# "syst*_u" (syst2, syst3, syst5 -> columns S1, T1, V1). The other syst*_c
# headers (syst0_c/Q1, syst1_c/R1, syst4_c/U1, syst6_c/W1) and syst_tot/N1
# are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("S1").Value = "syst2_u"
$ws.Range("T1").Value = "syst3_u"
$ws.Range("V1").Value = "syst5_u"

# Move the active selection to R16, matching the saved cursor position.
$ws.Range("R16").Select()
